# Update NATMI LR-pairs sheet (Efna2-Epha2) with newly-computed TPM numbers.
# The three "ECs" sending-cluster rows are removed entirely (only FAPs and
# MuSCs remain as sending clusters), and every remaining row's numeric
# columns (G:J, M:T) are refreshed with the new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the old rows whose "Sending cluster" (col A) was "ECs".
#    These were rows 2-4 (ECs->ECs, ECs->FAPs, ECs->MuSCs).
$ws.Rows("2:4").Delete()

# 2. Refresh the numeric columns for the remaining rows (now rows 2-7,
#    sending clusters FAPs then MuSCs, target clusters ECs/FAPs/MuSCs).

# Row 2: FAPs -> ECs
$ws.Range("G2").Value = 2.015377
$ws.Range("H2").Value = 6.046131
$ws.Range("I2").Value = 0.7554960962715589
$ws.Range("J2").Value = 0.7554960962715588
$ws.Range("M2").Value = 10.858287
$ws.Range("N2").Value = 32.574861
$ws.Range("O2").Value = 0.5084025289165609
$ws.Range("P2").Value = 0.508402528916561
$ws.Range("Q2").Value = 21.883541879199
$ws.Range("R2").Value = 196.951876912791
$ws.Range("S2").Value = 0.3840961259310501
$ws.Range("T2").Value = 0.3840961259310501

# Row 3: FAPs -> FAPs
$ws.Range("G3").Value = 2.015377
$ws.Range("H3").Value = 6.046131
$ws.Range("I3").Value = 0.7554960962715589
$ws.Range("J3").Value = 0.7554960962715588
$ws.Range("O3").Value = 0.004437346842596906
$ws.Range("P3").Value = 0.004437346842596906
$ws.Range("Q3").Value = 0.1909999654593333
$ws.Range("R3").Value = 1.718999689134
$ws.Range("S3").Value = 0.00335239821738489
$ws.Range("T3").Value = 0.00335239821738489

# Row 4: FAPs -> MuSCs
$ws.Range("G4").Value = 2.015377
$ws.Range("H4").Value = 6.046131
$ws.Range("I4").Value = 0.7554960962715589
$ws.Range("J4").Value = 0.7554960962715588
$ws.Range("O4").Value = 0.4871601242408422
$ws.Range("P4").Value = 0.4871601242408422
$ws.Range("Q4").Value = 20.969189518823
$ws.Range("R4").Value = 188.722705669407
$ws.Range("S4").Value = 0.3680475721231239
$ws.Range("T4").Value = 0.3680475721231238

# Row 5: MuSCs -> ECs
$ws.Range("G5").Value = 0.6522436666666667
$ws.Range("H5").Value = 1.956731
$ws.Range("I5").Value = 0.2445039037284412
$ws.Range("J5").Value = 0.2445039037284411
$ws.Range("M5").Value = 10.858287
$ws.Range("N5").Value = 32.574861
$ws.Range("O5").Value = 0.5084025289165609
$ws.Range("P5").Value = 0.508402528916561
$ws.Range("Q5").Value = 7.082248926598999
$ws.Range("R5").Value = 63.74024033939099
$ws.Range("S5").Value = 0.1243064029855108
$ws.Range("T5").Value = 0.1243064029855108

# Row 6: MuSCs -> FAPs
$ws.Range("G6").Value = 0.6522436666666667
$ws.Range("H6").Value = 1.956731
$ws.Range("I6").Value = 0.2445039037284412
$ws.Range("J6").Value = 0.2445039037284411
$ws.Range("O6").Value = 0.004437346842596906
$ws.Range("P6").Value = 0.004437346842596906
$ws.Range("Q6").Value = 0.06181400194822221
$ws.Range("R6").Value = 0.5563260175339999
$ws.Range("S6").Value = 0.001084948625212016
$ws.Range("T6").Value = 0.001084948625212016

# Row 7: MuSCs -> MuSCs
$ws.Range("G7").Value = 0.6522436666666667
$ws.Range("H7").Value = 1.956731
$ws.Range("I7").Value = 0.2445039037284412
$ws.Range("J7").Value = 0.2445039037284411
$ws.Range("O7").Value = 0.4871601242408422
$ws.Range("P7").Value = 0.4871601242408422
$ws.Range("Q7").Value = 6.786333801956333
$ws.Range("R7").Value = 61.077004217607
$ws.Range("S7").Value = 0.1191125521177183
$ws.Range("T7").Value = 0.1191125521177183
